{"js": "// Clarify the jetson io file name: \"io-config\" -> \"jetson-io\" in the two\n// paragraphs that reference the opt/nvidia config script and the\n// `sudo python3 ...` call that invokes it (but NOT the unrelated\n// \"io-config-by-pins.py\" verification script a few lines down).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text || \"\";\n  // Only touch the two specific lines; skip \"io-config-by-pins.py\".\n  if (text.indexOf(\"In opt/nvidia/io-config\") !== -1 ||\n      text.indexOf(\"Call sudo python3 io-config.py\") !== -1) {\n    const hits = paragraph.search(\"io-config\", { matchCase: true, matchWholeWord: false });\n    hits.load(\"items\");\n    await context.sync();\n    for (const hit of hits.items) {\n      hit.insertText(\"jetson-io\", \"Replace\");\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Clarify the jetson io file name: \"io-config\" -> \"jetson-io\" in the two\n# paragraphs that reference the opt/nvidia config script and the\n# `sudo python3 ...` call that invokes it (but NOT the unrelated\n# \"io-config-by-pins.py\" verification script a few lines down).\n$d = $word.ActiveDocument\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"In opt/nvidia/io-config*\" -or $t -like \"Call sudo python3 io-config.py*\") {\n        $r = $p.Range\n        [void]$r.Find.Execute(\"io-config\", $false, $false, $false, $false, $false, $true, 1, $false, \"jetson-io\", 2)\n    }\n}\n"}
